# Update countries & provincias Spain
#
# Daily COVID-19 data refresh: updates "Casos totales", "Nuevos casos",
# "Casos activos", "Recuperados", "Casos criticos", "Muertes hoy" and
# "Muertes" for the countries whose figures changed, then re-sorts the
# A4:H202 table by "Casos totales" (column B) descending -- some rows
# change position because their totals now outrank their neighbours.
# Also bumps the "last refreshed" timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Cells.Item(6, 2).Value = 54974
$ws.Cells.Item(6, 3).Value = 118
$ws.Cells.Item(6, 5).Value = 53810
$ws.Cells.Item(6, 7).Value = 5
$ws.Cells.Item(6, 8).Value = 785
# Row 13
$ws.Cells.Item(13, 2).Value = 8227
$ws.Cells.Item(13, 3).Value = 150
$ws.Cells.Item(13, 5).Value = 7665
$ws.Cells.Item(13, 7).Value = 5
$ws.Cells.Item(13, 8).Value = 427
# Row 17
$ws.Cells.Item(17, 1).Value = 'Portugal'
$ws.Cells.Item(17, 2).Value = 2995
$ws.Cells.Item(17, 3).Value = 633
$ws.Cells.Item(17, 4).Value = 22
$ws.Cells.Item(17, 5).Value = 2930
$ws.Cells.Item(17, 6).Value = 61
$ws.Cells.Item(17, 7).Value = 10
$ws.Cells.Item(17, 8).Value = 43
# Row 18
$ws.Cells.Item(18, 1).Value = 'Noruega'
$ws.Cells.Item(18, 2).Value = 2902
$ws.Cells.Item(18, 3).Value = 36
$ws.Cells.Item(18, 4).Value = 6
$ws.Cells.Item(18, 5).Value = 2883
$ws.Cells.Item(18, 6).Value = 57
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 8).Value = 13
# Row 19
$ws.Cells.Item(19, 1).Value = 'Canada'
$ws.Cells.Item(19, 2).Value = 2792
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 112
$ws.Cells.Item(19, 5).Value = 2654
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 26
# Row 20
$ws.Cells.Item(20, 1).Value = 'Australia'
$ws.Cells.Item(20, 2).Value = 2431
$ws.Cells.Item(20, 3).Value = 114
$ws.Cells.Item(20, 4).Value = 118
$ws.Cells.Item(20, 5).Value = 2304
$ws.Cells.Item(20, 6).Value = 11
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(20, 8).Value = 9
# Row 21
$ws.Cells.Item(21, 2).Value = 2345
$ws.Cells.Item(21, 3).Value = 46
$ws.Cells.Item(21, 5).Value = 2288
$ws.Cells.Item(21, 6).Value = 158
# Row 33
$ws.Cells.Item(33, 1).Value = 'Polonia'
$ws.Cells.Item(33, 2).Value = 957
$ws.Cells.Item(33, 3).Value = 56
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 5).Value = 943
$ws.Cells.Item(33, 6).Value = 3
$ws.Cells.Item(33, 7).Value = 3
$ws.Cells.Item(33, 8).Value = 13
# Row 34
$ws.Cells.Item(34, 1).Value = 'Tailandia'
$ws.Cells.Item(34, 2).Value = 934
$ws.Cells.Item(34, 3).Value = 107
$ws.Cells.Item(34, 4).Value = 70
$ws.Cells.Item(34, 5).Value = 860
$ws.Cells.Item(34, 6).Value = 11
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 4
# Row 37
$ws.Cells.Item(37, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(37, 2).Value = 900
$ws.Cells.Item(37, 3).Value = 133
$ws.Cells.Item(37, 4).Value = 28
$ws.Cells.Item(37, 5).Value = 871
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 1
# Row 38
$ws.Cells.Item(38, 1).Value = 'Finlandia'
$ws.Cells.Item(38, 2).Value = 853
$ws.Cells.Item(38, 3).Value = 61
$ws.Cells.Item(38, 4).Value = 10
$ws.Cells.Item(38, 5).Value = 840
$ws.Cells.Item(38, 6).Value = 11
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = 3
# Row 39
$ws.Cells.Item(39, 1).Value = 'Indonesia'
$ws.Cells.Item(39, 2).Value = 790
$ws.Cells.Item(39, 3).Value = 104
$ws.Cells.Item(39, 4).Value = 31
$ws.Cells.Item(39, 5).Value = 701
$ws.Cells.Item(39, 7).Value = 3
$ws.Cells.Item(39, 8).Value = 58
# Row 84
$ws.Cells.Item(84, 2).Value = 141
$ws.Cells.Item(84, 3).Value = 7
$ws.Cells.Item(84, 5).Value = 124
# Row 111
$ws.Cells.Item(111, 2).Value = 56
$ws.Cells.Item(111, 3).Value = 6
$ws.Cells.Item(111, 5).Value = 56
# Row 114
$ws.Cells.Item(114, 1).Value = 'Mauricio'
$ws.Cells.Item(114, 3).Value = 6
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 6).Value = 1
$ws.Cells.Item(114, 8).Value = 2
# Row 115
$ws.Cells.Item(115, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(115, 3).Value = 3
$ws.Cells.Item(115, 6).Value = 0
# Row 116
$ws.Cells.Item(116, 1).Value = 'Cuba'
$ws.Cells.Item(116, 2).Value = 48
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 1
$ws.Cells.Item(116, 5).Value = 46
$ws.Cells.Item(116, 6).Value = 2
# Row 117
$ws.Cells.Item(117, 1).Value = 'Nigeria'
$ws.Cells.Item(117, 2).Value = 46
$ws.Cells.Item(117, 3).Value = 2
$ws.Cells.Item(117, 4).Value = 2
$ws.Cells.Item(117, 8).Value = 1
# Row 123
$ws.Cells.Item(123, 1).Value = 'Honduras'
$ws.Cells.Item(123, 3).Value = 6
# Row 124
$ws.Cells.Item(124, 1).Value = 'Mayotte'
$ws.Cells.Item(124, 3).Value = 0
# Row 132
$ws.Cells.Item(132, 1).Value = 'Monaco'
$ws.Cells.Item(132, 3).Value = 0
# Row 133
$ws.Cells.Item(133, 1).Value = 'Togo'
$ws.Cells.Item(133, 3).Value = 3
# Row 141
$ws.Cells.Item(141, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(141, 3).Value = 4
# Row 142
$ws.Cells.Item(142, 1).Value = 'Uganda'
$ws.Cells.Item(142, 3).Value = 5
# Row 144
$ws.Cells.Item(144, 1).Value = 'Tanzania'
# Row 146
$ws.Cells.Item(146, 1).Value = 'Etiopia'
# Row 152
$ws.Cells.Item(152, 1).Value = 'Seychelles'
# Row 153
$ws.Cells.Item(153, 1).Value = 'Dominica'
# Row 154
$ws.Cells.Item(154, 1).Value = 'Surinam'
# Row 155
$ws.Cells.Item(155, 1).Value = 'Haiti'
# Row 158
$ws.Cells.Item(158, 1).Value = 'Benin'
# Row 159
$ws.Cells.Item(159, 1).Value = 'Bermudas'
# Row 160
$ws.Cells.Item(160, 1).Value = 'Gabon'
# Row 161
$ws.Cells.Item(161, 1).Value = 'Islas Caimanes'
# Row 164
$ws.Cells.Item(164, 1).Value = 'Bahamas'
$ws.Cells.Item(164, 4).Value = 1
$ws.Cells.Item(164, 8).Value = 0
# Row 165
$ws.Cells.Item(165, 1).Value = 'Guyana'
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 8).Value = 1
# Row 172
$ws.Cells.Item(172, 1).Value = 'Santa Lucia'
# Row 173
$ws.Cells.Item(173, 1).Value = 'Republica del Chad'
# Row 174
$ws.Cells.Item(174, 1).Value = 'Mozambique'
# Row 176
$ws.Cells.Item(176, 1).Value = 'Liberia'
# Row 177
$ws.Cells.Item(177, 1).Value = 'Republica de Africa Central'
# Row 178
$ws.Cells.Item(178, 1).Value = 'Laos'
$ws.Cells.Item(178, 3).Value = 1
# Row 179
$ws.Cells.Item(179, 1).Value = 'Antigua y Barbuda'
# Row 180
$ws.Cells.Item(180, 1).Value = 'Birmania'
$ws.Cells.Item(180, 3).Value = 0
# Row 181
$ws.Cells.Item(181, 1).Value = 'San Bartolome'
# Row 182
$ws.Cells.Item(182, 1).Value = 'Nepal'
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 1
$ws.Cells.Item(182, 8).Value = 0
# Row 183
$ws.Cells.Item(183, 1).Value = 'Zimbabue'
# Row 184
$ws.Cells.Item(184, 1).Value = 'Gambia'
# Row 185
$ws.Cells.Item(185, 1).Value = 'Sudan'
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 8).Value = 1
# Row 186
$ws.Cells.Item(186, 1).Value = 'Butan'
# Row 187
$ws.Cells.Item(187, 1).Value = 'San Martin (Parte Holandesa)'
# Row 188
$ws.Cells.Item(188, 1).Value = 'Nicaragua'
# Row 189
$ws.Cells.Item(189, 1).Value = 'Mauritania'
$ws.Cells.Item(189, 3).Value = 0
# Row 190
$ws.Cells.Item(190, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(190, 3).Value = 2
# Row 193
$ws.Cells.Item(193, 1).Value = 'Eritrea'
# Row 194
$ws.Cells.Item(194, 1).Value = 'Timor Oriental'
# Row 197
$ws.Cells.Item(197, 1).Value = 'Papua Nueva Guinea'
# Row 198
$ws.Cells.Item(198, 1).Value = 'Montserrat'
# Row 200
$ws.Cells.Item(200, 1).Value = 'San Vicente y las Granadinas'

# Update the "last refreshed" timestamp banner
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 25 de Marzo de 2020 a las 13:46"
